$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update URL (B2) and Date (B8) ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "https://interop.esante.gouv.fr/ig/fhir/ruim/ValueSet/MeltingPotVS"
$wsMeta.Range("B8").Value = "2026-02-06T10:39:27+00:00"

# --- Sheet "Include #1": update competence-code-system URL (B4) ---
$wsInc1 = $wb.Worksheets.Item("Include #1")
$wsInc1.Range("B4").Value = "https://interop.esante.gouv.fr/ig/fhir/ruim/CodeSystem/competence-code-system"

# --- Sheet "Include #2": update type-carte-code-system URL (B4) ---
$wsInc2 = $wb.Worksheets.Item("Include #2")
$wsInc2.Range("B4").Value = "https://interop.esante.gouv.fr/ig/fhir/ruim/CodeSystem/type-carte-code-system"
